$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.382.04"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.433.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "565.01"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.87"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.45%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.427.50"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.108"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.82%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.352"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.08%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -3.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000174"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.870.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.213.53"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.435.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.45"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.03"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.35%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.02"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "623.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.97"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0965"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -8.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.552.33"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.46"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.06"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.88%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.88%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.01"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.46"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.375"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.70"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.25"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -5.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "146.36"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.74"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -6.77%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.54"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.67%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "144.87"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.70"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.32"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -4.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0524"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.593"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0229"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.79%  "
